$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 128, shifting existing rows 128-149 down to 129-150
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new record
$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value = "La Araucanía"
$ws.Cells.Item(128, 4).Value = 44504
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = 100112005
$ws.Cells.Item(128, 7).Value = "Puerro"
$ws.Cells.Item(128, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 110
$ws.Cells.Item(128, 11).Value = 7000
$ws.Cells.Item(128, 12).Value = 7000
$ws.Cells.Item(128, 13).Value = 7000
$ws.Cells.Item(128, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(128, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(128, 16).Value = 583
$ws.Cells.Item(128, 17).Value = 12
$ws.Cells.Item(128, 18).Value = "Hortaliza"

# Apply the same number format style as other date cells in column D to the new D128
$ws.Cells.Item(128, 4).NumberFormat = $ws.Cells.Item(129, 4).NumberFormat
